$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.541.50'
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").Value = '3.562.63'
$ws.Range("E3").Value = '  +7.42%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''239.93'
$ws.Range("E5").Value = '  +4.51%  '
$ws.Range("D6").Value = '''637.32'
$ws.Range("E6").Value = '  +3.37%  '
$ws.Range("D7").Value = '''1.48'
$ws.Range("E7").Value = '  +8.94%  '
$ws.Range("D8").Value = '''0.400'
$ws.Range("E8").Value = '  +3.81%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +11.82%  '
$ws.Range("D11").Value = '3.557.66'
$ws.Range("E11").Value = '  +7.29%  '
$ws.Range("D12").Value = '''43.54'
$ws.Range("E12").Value = '  +4.35%  '
$ws.Range("D13").Value = '''0.201'
$ws.Range("E13").Value = '  +5.06%  '
$ws.Range("D14").Value = '''6.40'
$ws.Range("E14").Value = '  +8.11%  '
$ws.Range("D15").Value = '4.244.36'
$ws.Range("E15").Value = '  +7.95%  '
$ws.Range("D16").Value = '95.399.73'
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("E17").Value = '  +4.67%  '
$ws.Range("D18").Value = '''8.34'
$ws.Range("E18").Value = '  +4.13%  '
$ws.Range("D19").Value = '3.567.15'
$ws.Range("E19").Value = '  +7.63%  '
$ws.Range("D20").Value = '''12.89'
$ws.Range("E20").Value = '  +18.33%  '
$ws.Range("D21").Value = '''18.05'
$ws.Range("E21").Value = '  +5.48%  '
$ws.Range("D22").Value = '''0.508'
$ws.Range("E22").Value = '  +12.31%  '
$ws.Range("D23").Value = '''516.00'
$ws.Range("E23").Value = '  +5.33%  '
$ws.Range("D24").Value = '''3.43'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = '''6.80'
$ws.Range("E25").Value = '  +12.77%  '
$ws.Range("D26").Value = '''0.0000194'
$ws.Range("E26").Value = '  +7.77%  '
$ws.Range("D27").Value = '''93.04'
$ws.Range("E27").Value = '  +4.06%  '
$ws.Range("D28").Value = '''12.36'
$ws.Range("E28").Value = '  +6.91%  '
$ws.Range("D29").Value = '''3.03'
$ws.Range("E29").Value = '  +15.91%  '
$ws.Range("D30").Value = '''0.144'
$ws.Range("E30").Value = '  +5.70%  '
$ws.Range("D31").Value = '''11.54'
$ws.Range("E31").Value = '  +4.55%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = '''0.183'
$ws.Range("E33").Value = '  +7.03%  '
$ws.Range("D34").Value = '''1.01'
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").Value = '''30.24'
$ws.Range("E35").Value = '  +7.00%  '
$ws.Range("D36").Value = '''0.568'
$ws.Range("E36").Value = '  +8.04%  '
$ws.Range("D37").Value = '''584.71'
$ws.Range("E37").Value = '  +11.57%  '
$ws.Range("D38").Value = '''7.69'
$ws.Range("E38").Value = '  +5.81%  '
$ws.Range("E39").Value = '  +7.66%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''0.936'
$ws.Range("E40").Value = '  +8.38%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.151'
$ws.Range("E41").Value = '  +3.79%  '
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("E43").Value = '  +4.87%  '
$ws.Range("D44").Value = '''0.0428'
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("E46").Value = '  +5.62%  '
$ws.Range("D47").Value = '''3.59'
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("D48").Value = '''2.18'
$ws.Range("E48").Value = '  +5.01%  '
$ws.Range("D49").Value = '''54.03'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").Value = '''8.18'
$ws.Range("E50").Value = '  +3.56%  '
$ws.Range("D51").Value = '''3.13'
$ws.Range("E51").Value = '  +4.18%  '
